$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Förändrad") holds a date serial that was bumped by one day
# (2024-05-02 -> 2024-05-03, i.e. serial 45414 -> 45415) for rows 2 through 28.
for ($row = 2; $row -le 28; $row++) {
    $ws.Cells.Item($row, 3).Value = 45415
}
